# Insert a new row before row 320, shifting existing rows 320:439 down to 321:440.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("320:320").Insert()

# Populate the newly inserted row 320 with the new data entry.
$ws.Range("A320").Value = 1
$ws.Range("B320").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C320").Value = "Arica y Parinacota"
$ws.Range("D320").Value = 44825
$ws.Range("E320").Value = 15
$ws.Range("F320").Value = 100112023
$ws.Range("G320").Value = "Brócoli"
$ws.Range("H320").Value = "Sin especificar"
$ws.Range("I320").Value = "Tercera"
$ws.Range("J320").Value = 800
$ws.Range("K320").Value = 400
$ws.Range("L320").Value = 500
$ws.Range("M320").Value = 450
$ws.Range("N320").Value = "$/unidad"
$ws.Range("O320").Value = "Región de Arica y Parinacota"
$ws.Range("P320").Value = 450
$ws.Range("Q320").Value = 1
$ws.Range("R320").Value = "Hortaliza"
